$d = $word.ActiveDocument
$d.Content.Find.Execute("9126-3772", $true, $false, $false, $false, $false, $true, 1, $false, "99126-3772", 2)
